$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 4408.8335
$ws.Range("J112").Value = 4550.5293
$ws.Range("L112").Value = 13651.5879
$ws.Range("N112").Value = -15867.5879
$ws.Range("H137").Value = 1713.3077
$ws.Range("I137").Value = 1462.2142
$ws.Range("K137").Value = 4386.642599999999
$ws.Range("M137").Value = -1836.642599999999
$ws.Range("H138").Value = 3160.975
$ws.Range("I138").Value = 8864.5
$ws.Range("J138").Value = 2154.4707
$ws.Range("K138").Value = 26593.5
$ws.Range("L138").Value = 6463.4121
$ws.Range("M138").Value = -21453.5
$ws.Range("N138").Value = -16743.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2434.977
$ws.Range("I32").Value = 1656.0548
$ws.Range("K32").Value = 1656.0548
$ws.Range("M32").Value = -1369.0548
$ws.Range("H45").Value = 1572.5834
$ws.Range("I45").Value = 1266.3334
$ws.Range("K45").Value = 1266.3334
$ws.Range("M45").Value = -889.3334
$ws.Range("H61").Value = 2555.7083
$ws.Range("I61").Value = 1324.5
$ws.Range("J61").Value = 5018.125
$ws.Range("K61").Value = 1324.5
$ws.Range("L61").Value = 5018.125
$ws.Range("M61").Value = -1112.5
$ws.Range("N61").Value = -5442.125
$ws.Range("H63").Value = 8000.6665
$ws.Range("I63").Value = 7800.8
$ws.Range("K63").Value = 7800.8
$ws.Range("M63").Value = -7114.8
$ws.Range("H66").Value = 8000.6665
$ws.Range("I66").Value = 7800.8
$ws.Range("K66").Value = 39004
$ws.Range("M66").Value = -35572
$ws.Range("H74").Value = 1679.2106
$ws.Range("I74").Value = 1599.5555
$ws.Range("J74").Value = 1750.9
$ws.Range("K74").Value = 1599.5555
$ws.Range("L74").Value = 1750.9
$ws.Range("M74").Value = -725.5554999999999
$ws.Range("N74").Value = -3498.9
$ws.Range("H77").Value = 1679.2106
$ws.Range("I77").Value = 1599.5555
$ws.Range("J77").Value = 1750.9
$ws.Range("K77").Value = 7997.7775
$ws.Range("L77").Value = 8754.5
$ws.Range("M77").Value = -3629.7775
$ws.Range("N77").Value = -17490.5
$ws.Range("H132").Value = 1433.24
$ws.Range("I132").Value = 1067.1351
$ws.Range("K132").Value = 3201.4053
$ws.Range("M132").Value = -671.4052999999999
$ws.Range("H136").Value = 2555.7083
$ws.Range("I136").Value = 1324.5
$ws.Range("J136").Value = 5018.125
$ws.Range("K136").Value = 3973.5
$ws.Range("L136").Value = 15054.375
$ws.Range("M136").Value = -1423.5
$ws.Range("N136").Value = -20154.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H61").Value = 16000
$ws.Range("J61").Value = 16000
$ws.Range("L61").Value = 16000
$ws.Range("N61").Value = -16626
$ws.Range("H75").Value = 9933
$ws.Range("I75").Value = 9288.333000000001
$ws.Range("J75").Value = 11867
$ws.Range("K75").Value = 9288.333000000001
$ws.Range("L75").Value = 11867
$ws.Range("M75").Value = -8352.333000000001
$ws.Range("N75").Value = -13739
$ws.Range("H78").Value = 9933
$ws.Range("I78").Value = 9288.333000000001
$ws.Range("J78").Value = 11867
$ws.Range("K78").Value = 27864.999
$ws.Range("L78").Value = 35601
$ws.Range("M78").Value = -23184.999
$ws.Range("N78").Value = -44961
$ws.Range("H105").Value = 2367.348
$ws.Range("I105").Value = 2357.1667
$ws.Range("J105").Value = 2404
$ws.Range("K105").Value = 2357.1667
$ws.Range("L105").Value = 2404
$ws.Range("M105").Value = -610.1667000000002
$ws.Range("N105").Value = -5898

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1631.4062
$ws.Range("I58").Value = 958.5
$ws.Range("J58").Value = 3111.8
$ws.Range("K58").Value = 958.5
$ws.Range("L58").Value = 3111.8
$ws.Range("M58").Value = -755.5
$ws.Range("N58").Value = -3517.8
$ws.Range("H132").Value = 2526.1428
$ws.Range("I132").Value = 1594.6666
$ws.Range("J132").Value = 3600.923
$ws.Range("K132").Value = 4783.9998
$ws.Range("L132").Value = 10802.769
$ws.Range("M132").Value = -2253.9998
$ws.Range("N132").Value = -15862.769
$ws.Range("H134").Value = 1915.7142
$ws.Range("I134").Value = 1637.64
$ws.Range("K134").Value = 4912.92
$ws.Range("M134").Value = -2377.92
$ws.Range("H136").Value = 1631.4062
$ws.Range("I136").Value = 958.5
$ws.Range("J136").Value = 3111.8
$ws.Range("K136").Value = 2875.5
$ws.Range("L136").Value = 9335.400000000001
$ws.Range("M136").Value = -325.5
$ws.Range("N136").Value = -14435.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 569
$ws.Range("I26").Value = 569.6667
$ws.Range("J26").Value = 568.3333
$ws.Range("K26").Value = 1709.0001
$ws.Range("L26").Value = 1704.9999
$ws.Range("M26").Value = -1421.0001
$ws.Range("N26").Value = -2280.9999
$ws.Range("H104").Value = 3258.9524
$ws.Range("J104").Value = 3560.7222
$ws.Range("L104").Value = 10682.1666
$ws.Range("N104").Value = -15924.1666
$ws.Range("H121").Value = 795.1429000000001
$ws.Range("I121").Value = 533
$ws.Range("K121").Value = 1599
$ws.Range("M121").Value = -289
$ws.Range("H131").Value = 775.5
$ws.Range("I131").Value = 521
$ws.Range("J131").Value = 791.7447
$ws.Range("K131").Value = 1563
$ws.Range("L131").Value = 2375.2341
$ws.Range("M131").Value = 3477
$ws.Range("N131").Value = -12455.2341

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1023.25
$ws.Range("I80").Value = 997.6667
$ws.Range("J80").Value = 1100
$ws.Range("K80").Value = 997.6667
$ws.Range("L80").Value = 1100
$ws.Range("M80").Value = 0.3333000000000084
$ws.Range("N80").Value = -3096
$ws.Range("H83").Value = 1023.25
$ws.Range("I83").Value = 997.6667
$ws.Range("J83").Value = 1100
$ws.Range("K83").Value = 4988.3335
$ws.Range("L83").Value = 5500
$ws.Range("M83").Value = 3.666500000000269
$ws.Range("N83").Value = -15484
$ws.Range("H97").Value = 2010.1
$ws.Range("I97").Value = 2402
$ws.Range("K97").Value = 2402
$ws.Range("M97").Value = -1906
$ws.Range("H126").Value = 2755.037
$ws.Range("I126").Value = 2748.8845
$ws.Range("J126").Value = 2915
$ws.Range("K126").Value = 8246.6535
$ws.Range("L126").Value = 8745
$ws.Range("M126").Value = -5776.6535
$ws.Range("N126").Value = -13685
$ws.Range("H132").Value = 2544.4285
$ws.Range("I132").Value = 2302.7083
$ws.Range("K132").Value = 6908.124899999999
$ws.Range("M132").Value = -4378.124899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3273
$ws.Range("I7").Value = 2285.3
$ws.Range("K7").Value = 2285.3
$ws.Range("M7").Value = -2173.3
$ws.Range("H40").Value = 4179.222
$ws.Range("I40").Value = 2046.909
$ws.Range("J40").Value = 7530
$ws.Range("K40").Value = 2046.909
$ws.Range("L40").Value = 7530
$ws.Range("M40").Value = -1910.909
$ws.Range("N40").Value = -7802
$ws.Range("H122").Value = 5742.3335
$ws.Range("I122").Value = 3878.6667
$ws.Range("J122").Value = 11333.333
$ws.Range("K122").Value = 11636.0001
$ws.Range("L122").Value = 33999.999
$ws.Range("M122").Value = -9186.000100000001
$ws.Range("N122").Value = -38899.999
$ws.Range("H126").Value = 3273
$ws.Range("I126").Value = 2285.3
$ws.Range("K126").Value = 6855.900000000001
$ws.Range("M126").Value = -4385.900000000001
$ws.Range("H132").Value = 2537.4688
$ws.Range("I132").Value = 2252.25
$ws.Range("J132").Value = 2708.6
$ws.Range("K132").Value = 6756.75
$ws.Range("L132").Value = 8125.799999999999
$ws.Range("M132").Value = -4226.75
$ws.Range("N132").Value = -13185.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 36328.316
$ws.Range("I122").Value = 37991.57
$ws.Range("K122").Value = 113974.71
$ws.Range("M122").Value = -111524.71
